# Naresh Mail Id was added 30/08/2022 3.24 pm
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting")

# Row 3: To_MailSend -> append naresh.kumar@e5.ai at the end
$ws.Range("B3").Value = " mvprasanth97@gmail.com ; lakshmi.u@tiliconveli.com ; narenbagavathye5@gmail.com ; sornalakshmie5@gmail.com ; aartiak.e5@gmail.com ; sankaravenie5@gmail.com ;sharongiftae5@gmail.com;naresh.kumar@e5.ai"

# Row 13: remove alagappan.m@e5.ai, keep only naresh.kumar@e5.ai
$ws.Range("B13").Value = "naresh.kumar@e5.ai"

# Row 14: remove alagappan.m@e5.ai from the list
$ws.Range("B14").Value = "naresh.kumar@e5.ai ; mvprasanth97@gmail.com ; lakshmi.u@tiliconveli.com ; narenbagavathye5@gmail.com ; sornalakshmie5@gmail.com ; aartiak.e5@gmail.com ; sankaravenie5@gmail.com;sharongiftae5@gmail.com"

# Move the active selection to B16 like in the edited file
$ws.Range("B16").Select()
